# Update: po 07. 06. 2021
# Revises several historical AgTests/AgPosit (columns F/G) figures and
# appends three new daily rows (2021-06-04, 2021-06-05, 2021-06-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing rows (columns F = AgTests, G = AgPosit) ---
$updates = @(
    @{ Cell = "F309"; Value = 78041 },
    @{ Cell = "G309"; Value = 5530 },
    @{ Cell = "F322"; Value = 109742 },
    @{ Cell = "G322"; Value = 2325 },
    @{ Cell = "F324"; Value = 250185 },
    @{ Cell = "F325"; Value = 774532 },
    @{ Cell = "G325"; Value = 6501 },
    @{ Cell = "F328"; Value = 181057 },
    @{ Cell = "F329"; Value = 73333 },
    @{ Cell = "G329"; Value = 1723 },
    @{ Cell = "F332"; Value = 484790 },
    @{ Cell = "G332"; Value = 4806 },
    @{ Cell = "F336"; Value = 82035 },
    @{ Cell = "F337"; Value = 105092 },
    @{ Cell = "F338"; Value = 221319 },
    @{ Cell = "F353"; Value = 723957 },
    @{ Cell = "G353"; Value = 5292 },
    @{ Cell = "F357"; Value = 138367 },
    @{ Cell = "G357"; Value = 3011 },
    @{ Cell = "F360"; Value = 748576 },
    @{ Cell = "G360"; Value = 5133 },
    @{ Cell = "F374"; Value = 773295 },
    @{ Cell = "G374"; Value = 3419 },
    @{ Cell = "F377"; Value = 176613 },
    @{ Cell = "F388"; Value = 730214 },
    @{ Cell = "G388"; Value = 2201 },
    @{ Cell = "F398"; Value = 298795 },
    @{ Cell = "F417"; Value = 343007 },
    @{ Cell = "F426"; Value = 107119 },
    @{ Cell = "G431"; Value = 403 },
    @{ Cell = "F434"; Value = 78972 },
    @{ Cell = "F442"; Value = 70345 },
    @{ Cell = "F443"; Value = 106861 },
    @{ Cell = "F449"; Value = 59821 },
    @{ Cell = "G451"; Value = 115 },
    @{ Cell = "F454"; Value = 51669 },
    @{ Cell = "F455"; Value = 50506 },
    @{ Cell = "G455"; Value = 120 },
    @{ Cell = "F456"; Value = 49682 },
    @{ Cell = "G456"; Value = 138 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# --- Append new daily rows 457-459 ---
$newRows = @(
    @{ Row = 457; A = 44351; B = 390340; C = 5548; D = 117; E = 12395; F = 71847; G = 121 },
    @{ Row = 458; A = 44352; B = 390436; C = 2780; D = 96;  E = 12404; F = 62290; G = 71  },
    @{ Row = 459; A = 44353; B = 390451; C = 981;  D = 15;  E = 12414; F = 52608; G = 76  }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").NumberFormat = "yyyy-mm-dd"
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
}
